$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.338
$ws.Range("B3").Value = 5.461
$ws.Range("D3").Value = -7.623
$ws.Range("E6").Value = 16.73
$ws.Range("D12").Value = -7.347
$ws.Range("B14").Value = 5.857000000000001
$ws.Range("E19").Value = 16.547
$ws.Range("B21").Value = 9.643000000000001
$ws.Range("B23").Value = 8.004999999999999
$ws.Range("D24").Value = -7.118
$ws.Range("E24").Value = 16.906
$ws.Range("B25").Value = 5.540999999999999
$ws.Range("C25").Value = -12.964
$ws.Range("D25").Value = -8.306999999999999
$ws.Range("B26").Value = 6.051
$ws.Range("C27").Value = -13.253
$ws.Range("B29").Value = 5.656000000000001
$ws.Range("E30").Value = 16.485
$ws.Range("C31").Value = -12.435
$ws.Range("E31").Value = 16.835
$ws.Range("E33").Value = 17.413
$ws.Range("C39").Value = -12.847
$ws.Range("E42").Value = 16.608
$ws.Range("C48").Value = -11.25
$ws.Range("D50").Value = -8.270000000000001
$ws.Range("C51").Value = -11.701
$ws.Range("C52").Value = -11.221
$ws.Range("B53").Value = 5.775
$ws.Range("D53").Value = -7.311
$ws.Range("C55").Value = -13.879
$ws.Range("E55").Value = 16.409
$ws.Range("C56").Value = -12.753
$ws.Range("B57").Value = 5.23
$ws.Range("C57").Value = -13.691
$ws.Range("D57").Value = -8.226999999999999
$ws.Range("E58").Value = 16.838
$ws.Range("B59").Value = 4.714
$ws.Range("D61").Value = -7.822000000000001
$ws.Range("D63").Value = -7.984999999999999
$ws.Range("E65").Value = 17.395
$ws.Range("B69").Value = 5.372999999999999
$ws.Range("D70").Value = -7.531000000000001
$ws.Range("E70").Value = 17.54
$ws.Range("C73").Value = -12.654
$ws.Range("E75").Value = 16.486
$ws.Range("B79").Value = 5.771
$ws.Range("B83").Value = 5.542
$ws.Range("E83").Value = 16.527
$ws.Range("D86").Value = -8.321999999999999
$ws.Range("E86").Value = 16.695
$ws.Range("C89").Value = -12.552
$ws.Range("C90").Value = -12.737
$ws.Range("B91").Value = 5.862
$ws.Range("C92").Value = -11.559
$ws.Range("B93").Value = 5.673
$ws.Range("E96").Value = 16.454
$ws.Range("E97").Value = 17.058
$ws.Range("D98").Value = -8.486000000000001
$ws.Range("D100").Value = -8.409000000000001
$ws.Range("D102").Value = -7.958000000000001
